# Auto-generated: updates cryptos price/volume figures (and re-sorts three
# rows by rank) to match the target OOXML diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.093.91"

$ws.Range("D3").Value = "2.576.50"
$ws.Range("E3").Value = "  -2.27%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'535.00"
$ws.Range("E5").Value = "  -0.37%  "

$ws.Range("D6").Value = "'141.46"
$ws.Range("E6").Value = "  -1.91%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  +3.56%  "

$ws.Range("D9").Value = "'6.73"
$ws.Range("E9").Value = "  +2.78%  "

$ws.Range("D10").Value = "'0.0992"
$ws.Range("E10").Value = "  -3.48%  "

$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("E12").Value = "  -1.81%  "

$ws.Range("D13").Value = "3.030.22"
$ws.Range("E13").Value = "  -2.22%  "

$ws.Range("D14").Value = "58.030.14"
$ws.Range("E14").Value = "  -2.89%  "

$ws.Range("D15").Value = "'20.68"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "2.570.94"
$ws.Range("E16").Value = "  -3.44%  "

$ws.Range("E17").Value = "  -2.91%  "

$ws.Range("E18").Value = "  -0.65%  "

$ws.Range("D19").Value = "'333.78"
$ws.Range("E19").Value = "  -2.62%  "

$ws.Range("D20").Value = "'10.01"
$ws.Range("E20").Value = "  -1.94%  "

$ws.Range("D21").Value = "'6.15"
$ws.Range("E21").Value = "  -4.00%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").Value = "'66.65"
$ws.Range("E23").Value = "  -1.27%  "

$ws.Range("E24").Value = "  +1.76%  "

$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("D26").Value = "'0.157"
$ws.Range("E26").Value = "  -5.40%  "

$ws.Range("E27").Value = "  -3.38%  "

$ws.Range("E28").Value = "  +0.06%  "

$ws.Range("D29").Value = "0.0₃0722"
$ws.Range("E29").Value = "  -4.01%  "

$ws.Range("E30").Value = "  -2.34%  "

$ws.Range("D31").Value = "'155.43"
$ws.Range("E31").Value = "  +2.84%  "

$ws.Range("D32").Value = "'5.83"
$ws.Range("E32").Value = "  -0.52%  "

$ws.Range("D33").Value = "'18.80"
$ws.Range("E33").Value = "  -0.53%  "

$ws.Range("E34").Value = "  -3.54%  "

$ws.Range("D35").Value = "'36.89"
$ws.Range("E35").Value = "  -1.43%  "

$ws.Range("D37").Value = "'0.826"
$ws.Range("E37").Value = "  +0.22%  "

$ws.Range("D38").Value = "'0.819"
$ws.Range("E38").Value = "  -2.47%  "

$ws.Range("E39").Value = "  -3.67%  "

$ws.Range("D40").Value = "'3.59"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").Value = "'282.01"
$ws.Range("E41").Value = "  -3.22%  "

$ws.Range("D42").Value = "'0.999"
$ws.Range("E42").Value = "  +0.02%  "

$ws.Range("B43").Value = "WhiteBITCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D43").Value = "'10.65"
$ws.Range("E43").Value = "  -0.81%  "

$ws.Range("B44").Value = "Mantle"
$ws.Range("C44").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.588"
$ws.Range("E44").Value = "  -2.50%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "'0.0949"
$ws.Range("E45").Value = "  -0.24%  "

$ws.Range("D46").Value = "'0.0533"
$ws.Range("E46").Value = "  +0.23%  "

$ws.Range("D47").Value = "'18.28"
$ws.Range("E47").Value = "  -4.54%  "

$ws.Range("E48").Value = "  +0.36%  "

$ws.Range("D49").Value = "1.904.73"
$ws.Range("E49").Value = "  -3.07%  "

$ws.Range("D50").Value = "'17.74"
$ws.Range("E50").Value = "  -4.09%  "

$ws.Range("D51").Value = "'4.35"
$ws.Range("E51").Value = "  -4.40%  "
